$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1) Update missing/present values on rows that stay (apply BEFORE the
#        row deletions below, while row numbers still match the original
#        layout: RM 125=19, RM135=21, RM140=23, SC5=27, SC101=29, SC119=31,
#        SC232=35) ---

# RM 125 (row 19): F19 was missing -> now has a value
$ws.Range("F19").Value = 17.81

# RM 135 (row 21): F21 had a value -> now missing
$ws.Range("F21").ClearContents()

# RM 140 (row 23): F23 was missing -> now has a value
$ws.Range("F23").Value = 16.48

# SC 5 (row 27): E27 had a value -> now missing
$ws.Range("E27").ClearContents()

# SC 101 (row 29): E29 was missing -> now has a value; F29 had a value -> now missing
$ws.Range("E29").Value = -10
$ws.Range("F29").ClearContents()

# SC 119 (row 31): E31 had a value -> now missing
$ws.Range("E31").ClearContents()

# SC 232 (row 35): F35 was missing -> now has a value
$ws.Range("F35").Value = 17.53

# --- 2) Delete the two rows that were dropped entirely: "RM 232" (row 26)
#        and "SC 92" (which becomes row 27 after the first delete shifts
#        everything up) ---
$ws.Rows.Item(26).Delete()
$ws.Rows.Item(27).Delete()
